$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "additional options" prompt text in D3: the option list is
# renumbered (a new "1 = 옵션없음" / no-option choice is introduced, and the
# old "없으시다면 0을 입력해주세요" line is dropped), shifting every other
# option down by one and adding a new "8 = 샷 + 시럽 + 사이즈업" combo.
$newText = "추가 옵션을 선택해주세요
1 = 옵션없음
2 = 샷추가
3 = 시럽추가
4 = 사이즈업
5 = 샷 + 시럽추가
6 = 샷 + 사이즈업
7 = 시럽 + 사이즈업
8 = 샷 + 시럽 + 사이즈업"

$cell = $ws.Range("D3")
$cell.Value = $newText

# The cell picks up a dedicated font/style for this text (new font + cellXf
# entry in styles.xml) — applying a font name touches that.
$cell.Font.Name = "Arial"
$cell.Font.Size = 12

# Move the active selection to D5 (previously E3).
$ws.Range("D5").Select()
